# Updated symbol list on Fri Dec 16 11:22:33 UTC 2022 with GitHub Actions
#
# All "Price" (column D) values in this sheet are stored as TEXT (they were
# written as inline/shared strings, not numbers) even though they look like
# numbers. Plain `.Value = "251.54"` assignment would make Excel's COM layer
# auto-coerce that into a real number, which would not match the source
# workbook's cell type. Prefixing the literal with a leading single-quote
# (the normal Excel "force text" idiom) keeps it text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - BNB
$ws.Range("D2").Value = "'251.54"

# Row 3 - OKB
$ws.Range("D3").Value = "'23.84"

# Row 4 - HuobiToken
$ws.Range("D4").Value = "'6.019"

# Row 5 - Cronos
$ws.Range("D5").Value = "'0.05967"

# Row 6
$ws.Range("D6").Value = "'3.420"

# Row 7
$ws.Range("D7").Value = "'6.571"

# Row 8
$ws.Range("D8").Value = "'1.331"

# Row 9
$ws.Range("D9").Value = "'0.7965"

# Row 10
$ws.Range("D10").Value = "'0.1494"

# Row 11
$ws.Range("D11").Value = "'0.07907"

# Row 12
$ws.Range("D12").Value = "'0.03354"

# Row 14
$ws.Range("D14").Value = "'0.09274"

# Row 15
$ws.Range("D15").Value = "'3.569"

# Row 16
$ws.Range("D16").Value = "'0.001668"

# Row 17
$ws.Range("D17").Value = "'0.04777"

# Row 18 - One: price + Volume(1h) label (lost "Worstin24h" suffix)
$ws.Range("D18").Value = "'0.0006085"
$ws.Range("E18").Value = "17OneONE"

# Row 19
$ws.Range("D19").Value = "'0.006245"

# Row 20
$ws.Range("D20").Value = "'0.005687"

# Row 21
$ws.Range("D21").Value = "'0.001068"

# Row 24
$ws.Range("D24").Value = "'2.213"

# Row 26
$ws.Range("D26").Value = "'0.1256"

# Row 27
$ws.Range("D27").Value = "'0.0006480"

# Row 40
$ws.Range("D40").Value = "'0.04432"

# Row 41
$ws.Range("D41").Value = "'0.007060"

# Row 42 and 43 swapped places (CEJI <-> BKEXToken), each keeping their own
# updated price/link/volume-label values.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1070"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.003153"
$ws.Range("E43").Value = "42CEJICEJI"

# Row 45
$ws.Range("D45").Value = "'0.002462"

# Row 46
$ws.Range("D46").Value = "'0.00005891"

# Row 48
$ws.Range("D48").Value = "'0.9914"

# Row 49 - BOLO: price + Volume(1h) label (gained "Worstin24h" suffix)
$ws.Range("D49").Value = "'0.1150"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"
